$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solver")

# Clear C3, set new values per diff
$ws.Range("C3").Value = $null
$ws.Range("F3").Value = 1
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 2

$ws.Range("E4").Value = 5
$ws.Range("H4").Value = 3
$ws.Range("J4").Value = 9

$ws.Range("H5").Value = 2

$ws.Range("E6").Value = 6
$ws.Range("K6").Value = 1

$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 4
$ws.Range("I7").Value = 7

$ws.Range("C8").Value = 9
$ws.Range("I8").Value = 5

$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 9

$ws.Range("D10").Value = 5
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 7

$ws.Range("D11").Value = 1
$ws.Range("G11").Value = 8

# Update selection to match diff (active cell J5)
$ws.Activate()
$ws.Range("J5").Select()
